$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 108-110, pushing the existing 108-110 rows
# (and their formatting) down to 111-113.
$ws.Range("A108:H110").Insert()

# The Insert() above does not preserve the original cell formatting (border /
# number format) for the newly inserted blank rows, so copy the formatting
# from the (now shifted) rows immediately below into the new rows.
$ws.Range("A111:H111").Copy()
$ws.Range("A108:H108").PasteSpecial(-4122)
$ws.Range("A112:H112").Copy()
$ws.Range("A109:H109").PasteSpecial(-4122)
$ws.Range("A113:H113").Copy()
$ws.Range("A110:H110").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$dot = [char]0x00B7

# New row 108
$ws.Range("A108").Value2 = 212
$ws.Range("B108").Value2 = "New ListingNVIDIA RTX 3060 Ti GIGABYTE EAGLE OC 8GB IN HAND READY TO SHIP"
$ws.Range("C108").Value2 = "Brand New $dot GIGABYTE"
$ws.Range("D108").Value2 = 650
$ws.Range("E108").Value2 = 10.9
$ws.Range("F108").Value2 = 660.9
$ws.Range("G108").Value2 = 44167
$ws.Range("H108").Value2 = "https://www.ebay.com/itm/NVIDIA-RTX-3060-Ti-GIGABYTE-EAGLE-OC-8GB-IN-HAND-READY-TO-SHIP/333811666219?hash=item4db8b9cd2b:g:tpEAAOSwbLxfyADU"

# New row 109
$ws.Range("A109").Value2 = 213
$ws.Range("B109").Value2 = "New ListingNVIDIA GeForce RTX 3060 TI Founders Edition IN HAND!!!! Limited Stock BUY NOW!!!"
$ws.Range("C109").Value2 = "Brand New $dot Geforce $dot 8 GB"
$ws.Range("D109").Value2 = 450
$ws.Range("E109").Value2 = 0
$ws.Range("F109").Value2 = 450
$ws.Range("G109").Value2 = 44167
$ws.Range("H109").Value2 = "https://www.ebay.com/itm/NVIDIA-GeForce-RTX-3060-TI-Founders-Edition-IN-HAND-Limited-Stock-BUY-NOW/203206293280?hash=item2f5009f720:g:e7MAAOSwr15fx~hb"

# New row 110
$ws.Range("A110").Value2 = 214
$ws.Range("B110").Value2 = "New ListingNVIDIA GeForce RTX 3060 Ti 8GB GDDR6PCI Express 4.0 Graphics Card"
$ws.Range("C110").Value2 = "Brand New $dot Geforce $dot 8 GB"
$ws.Range("D110").Value2 = 399.99
$ws.Range("E110").Value2 = 0
$ws.Range("F110").Value2 = 399.99
$ws.Range("G110").Value2 = 44167
$ws.Range("H110").Value2 = "https://www.ebay.com/itm/NVIDIA-GeForce-RTX-3060-Ti-8GB-GDDR6PCI-Express-4-0-Graphics-Card/313326954368?hash=item48f3bde780:g:~wwAAOSww5hfx9YC"

# Rows 111-113 (previously 108-110) keep almost all of their original values;
# only the "A" id column changes for each of them.
$ws.Range("A111").Value2 = 250
$ws.Range("A112").Value2 = 310
$ws.Range("A113").Value2 = 317

Write-Output "done"
